# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the upstream data source.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1206
$ws1.Range("F4").Value = 14531
$ws1.Range("F5").Value = 17423
$ws1.Range("F9").Value = 209
$ws1.Range("F12").Value = 0
$ws1.Range("F13").Value = 1
$ws1.Range("F16").Value = 25
$ws1.Range("F19").Value = 1315
$ws1.Range("F23").Value = 141
$ws1.Range("F24").Value = 7163
$ws1.Range("F27").Value = 1165
$ws1.Range("F29").Value = 5842
$ws1.Range("F31").Value = 38
$ws1.Range("F35").Value = 5048

# Sheet "全部类型" (All types) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1206
$ws4.Range("F4").Value = 14531
$ws4.Range("F5").Value = 17423
$ws4.Range("F9").Value = 209
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 1
$ws4.Range("F16").Value = 25
$ws4.Range("F19").Value = 1315
$ws4.Range("F24").Value = 141
$ws4.Range("F25").Value = 7163
$ws4.Range("F28").Value = 1165
$ws4.Range("F31").Value = 5842
$ws4.Range("F33").Value = 38
$ws4.Range("F37").Value = 5048
